# This workbook contains a single weekly price-listing sheet for
# "Hortaliza, Agrícola del Norte S.A. de Arica - Brócoli" with data rows
# running from row 2 (oldest entries first is not guaranteed, but every row
# has identical metadata columns) down to row 306.
#
# The commit adds one new weekly record. It is inserted at row 204, which
# pushes all the existing records (rows 204-306) down by one row (to rows
# 205-307) and grows the sheet from 306 to 307 data/header rows overall.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 204:306 down to 205:307, duplicating formatting from row 204.
$ws.Rows("204:204").Insert()

# Populate the newly inserted row 204 with the new weekly record.
$ws.Range("A204").Value = 1
$ws.Range("B204").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C204").Value = "Arica y Parinacota"
$ws.Range("D204").Value = 44609
$ws.Range("E204").Value = 15
$ws.Range("F204").Value = 100112023
$ws.Range("G204").Value = "Brócoli"
$ws.Range("H204").Value = "Sin especificar"
$ws.Range("I204").Value = "Tercera"
$ws.Range("J204").Value = 1600
$ws.Range("K204").Value = 350
$ws.Range("L204").Value = 400
$ws.Range("M204").Value = 375
$ws.Range("N204").Value = "`$/unidad"
$ws.Range("O204").Value = "Región de Arica y Parinacota"
$ws.Range("P204").Value = 375
$ws.Range("Q204").Value = 1
$ws.Range("R204").Value = "Hortaliza"
